$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(0.8080066826315555, 0.6215487750042737, 0.782496592585249, 2.9777541793440783)"
$ws.Range("C2").Value = "NIG(1.5587776180831914, 1.1812542410549782, 3.291643814789671, 6.63697551405021)"
$ws.Range("D2").Value = "MIE(10.720906307688107, 5.921097231636256, -10.123067574833264, 12.348964034806233)"
$ws.Range("E2").Value = "EXN(3.8008014316855077, 2.696047900087758, 2.4920936521089403)"
